$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the legend / footer row (row 24) with the missing values
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = "-"
$ws.Range("F24").Value = "-"
$ws.Range("G24").Value = "-"
$ws.Range("H24").Value = 4

# Move the active selection to match the author's last selection
$ws.Range("H25").Select()
